$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, reusing the same formatting as the
# other header cells (bold, bordered, centered style already in G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values for the two data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
